$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 2767  # F3: 2760 -> 2767
$ws.Cells.Item(5, 6).Value = 945  # F5: 944 -> 945
$ws.Cells.Item(7, 6).Value = 2432  # F7: 2423 -> 2432
$ws.Cells.Item(8, 6).Value = 1874  # F8: 1870 -> 1874
$ws.Cells.Item(9, 6).Value = 231  # F9: 229 -> 231
$ws.Cells.Item(11, 6).Value = 2528  # F11: 2526 -> 2528
$ws.Cells.Item(12, 6).Value = 568  # F12: 567 -> 568
$ws.Cells.Item(13, 6).Value = 259  # F13: 258 -> 259
$ws.Cells.Item(14, 6).Value = 4  # F14: 2 -> 4
$ws.Cells.Item(16, 6).Value = 135  # F16: 134 -> 135
$ws.Cells.Item(18, 6).Value = 9431  # F18: 9412 -> 9431
$ws.Cells.Item(21, 6).Value = 7389  # F21: 7364 -> 7389
$ws.Cells.Item(22, 6).Value = 11935  # F22: 11912 -> 11935
$ws.Cells.Item(26, 6).Value = 375  # F26: 374 -> 375
$ws.Cells.Item(28, 6).Value = 2679  # F28: 2673 -> 2679
$ws.Cells.Item(30, 6).Value = 210  # F30: 209 -> 210
$ws.Cells.Item(31, 6).Value = 2653  # F31: 2647 -> 2653
$ws.Cells.Item(32, 6).Value = 943  # F32: 926 -> 943
$ws.Cells.Item(36, 6).Value = 4545  # F36: 4544 -> 4545
$ws.Cells.Item(37, 6).Value = 1042  # F37: 1040 -> 1042
$ws.Cells.Item(38, 6).Value = 25  # F38: 24 -> 25
$ws.Cells.Item(39, 6).Value = 359  # F39: 358 -> 359

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 12  # F11: 11 -> 12
$ws.Cells.Item(16, 6).Value = 18  # F16: 14 -> 18
$ws.Cells.Item(19, 6).Value = 5  # F19: 4 -> 5
$ws.Cells.Item(20, 6).Value = 101  # F20: 100 -> 101

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 633  # F2: 632 -> 633
$ws.Cells.Item(4, 6).Value = 176  # F4: 173 -> 176

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 633  # F2: 632 -> 633
$ws.Cells.Item(4, 6).Value = 176  # F4: 173 -> 176
$ws.Cells.Item(6, 6).Value = 2767  # F6: 2760 -> 2767
$ws.Cells.Item(9, 6).Value = 945  # F9: 944 -> 945
$ws.Cells.Item(12, 6).Value = 2432  # F12: 2423 -> 2432
$ws.Cells.Item(14, 6).Value = 1874  # F14: 1870 -> 1874
$ws.Cells.Item(15, 6).Value = 231  # F15: 229 -> 231
$ws.Cells.Item(16, 6).Value = 2528  # F16: 2526 -> 2528
$ws.Cells.Item(17, 6).Value = 568  # F17: 567 -> 568
$ws.Cells.Item(18, 6).Value = 259  # F18: 258 -> 259
$ws.Cells.Item(19, 6).Value = 135  # F19: 134 -> 135
$ws.Cells.Item(21, 6).Value = 9432  # F21: 9412 -> 9432
$ws.Cells.Item(23, 6).Value = 7389  # F23: 7365 -> 7389
$ws.Cells.Item(24, 6).Value = 11935  # F24: 11912 -> 11935
$ws.Cells.Item(25, 6).Value = 12  # F25: 11 -> 12
$ws.Cells.Item(28, 6).Value = 375  # F28: 374 -> 375
$ws.Cells.Item(32, 6).Value = 2679  # F32: 2673 -> 2679
$ws.Cells.Item(34, 6).Value = 18  # F34: 14 -> 18
$ws.Cells.Item(36, 6).Value = 210  # F36: 209 -> 210
$ws.Cells.Item(39, 6).Value = 4545  # F39: 4544 -> 4545
$ws.Cells.Item(41, 6).Value = 5  # F41: 4 -> 5
$ws.Cells.Item(42, 6).Value = 101  # F42: 100 -> 101

